$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old "per capita" labels to new "per cap." labels (fix typo in label)
$replacements = @{
    "ln(GDP [dollars per capita])" = "ln(GDP [dollars per cap.])"
    "ln(Tourism - Inbound [per capita])" = "ln(Tourism - Inbound [per cap.])"
    "ln(Migrant Population [per capita])" = "ln(Migrant Population [per cap.])"
    "ln(ProMed Mentions [per capita])" = "ln(ProMed Mentions [per cap.])"
    "ln(AB Exports [dollars per capita])" = "ln(AB Exports [dollars per cap.])"
    "ln(Publication Bias Index [per capita])" = "ln(Publication Bias Index [per cap.])"
    "Livestock AB Consumption [kg per capita)" = "Livestock AB Consumption (kg per cap.)"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
